$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain plain text so values such as
# "249.56" or "1.00" are not reinterpreted as numbers (matches the
# original inlineStr text cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '36.763.20'
$ws.Range('E2').Value = '  +4.38%  '
$ws.Range('D3').Value = '1.923.76'
$ws.Range('E3').Value = '  +2.59%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '249.56'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').Value = '0.698'
$ws.Range('E6').Value = '  +2.78%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '44.40'
$ws.Range('E8').Value = '  +2.62%  '
$ws.Range('D9').Value = '58.61'
$ws.Range('E9').Value = '  +9.55%  '
$ws.Range('D10').Value = '0.368'
$ws.Range('E10').Value = '  +4.39%  '
$ws.Range('E11').Value = '  +4.17%  '
$ws.Range('E12').Value = '  +2.95%  '
$ws.Range('D13').Value = '14.64'
$ws.Range('E13').Value = '  +9.35%  '
$ws.Range('D14').Value = '0.803'
$ws.Range('E14').Value = '  +5.15%  '
$ws.Range('D15').Value = '2.199.47'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('E16').Value = '  +5.18%  '
$ws.Range('D17').Value = '1.918.95'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = '36.700.89'
$ws.Range('E18').Value = '  +4.14%  '
$ws.Range('D19').Value = '74.25'
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('D20').Value = '0.0₃0860'
$ws.Range('E20').Value = '  +5.41%  '
$ws.Range('D21').Value = '251.93'
$ws.Range('E21').Value = '  +3.72%  '
$ws.Range('D22').Value = '13.30'
$ws.Range('E22').Value = '  +4.68%  '
$ws.Range('D23').Value = '5.21'
$ws.Range('E23').Value = '  +5.97%  '
$ws.Range('D24').Value = '2.67'
$ws.Range('E24').Value = '  +2.20%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('E26').Value = '  +2.76%  '
$ws.Range('D27').Value = '168.56'
$ws.Range('E27').Value = '  +2.01%  '
$ws.Range('D28').Value = '8.85'
$ws.Range('E28').Value = '  +4.80%  '
$ws.Range('D29').Value = '18.85'
$ws.Range('E29').Value = '  +3.79%  '
$ws.Range('E30').Value = '  +2.73%  '
$ws.Range('D31').Value = '4.56'
$ws.Range('E31').Value = '  +7.19%  '
$ws.Range('D32').Value = '0.0619'
$ws.Range('E32').Value = '  +5.37%  '
$ws.Range('D33').Value = '1.99'
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('D34').Value = '4.37'
$ws.Range('E34').Value = '  +6.08%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').Value = '0.0862'
$ws.Range('E36').Value = '  +18.25%  '
$ws.Range('D37').Value = '1.51'
$ws.Range('E37').Value = '  -10.64%  '
$ws.Range('D38').Value = '0.895'
$ws.Range('E38').Value = '  +7.69%  '
$ws.Range('D39').Value = '17.86'
$ws.Range('E39').Value = '  +50.40%  '
$ws.Range('E40').Value = '  +5.15%  '
$ws.Range('D41').Value = '105.95'
$ws.Range('E41').Value = '  +10.90%  '
$ws.Range('E42').Value = '  +5.94%  '
$ws.Range('D43').Value = '17.38'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('E44').Value = '  +4.24%  '
$ws.Range('D45').Value = '1.338.94'
$ws.Range('E45').Value = '  +3.06%  '
$ws.Range('D46').Value = '2.56'
$ws.Range('E46').Value = '  +7.77%  '
$ws.Range('D47').Value = '2.39'
$ws.Range('E47').Value = '  +2.12%  '
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('D49').Value = '2.80'
$ws.Range('E49').Value = '  +3.14%  '
$ws.Range('D50').Value = '6.45'
$ws.Range('E50').Value = '  +4.23%  '
$ws.Range('D51').Value = '43.47'
$ws.Range('E51').Value = '  +4.38%  '
